# Update gh-pages to output generated at 456a3b4
$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 7730
$ws1.Range("F17").Value = 5585
$ws1.Range("F18").Value = 150
$ws1.Range("F20").Value = 995
$ws1.Range("F22").Value = 322

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 7730
$ws4.Range("F18").Value = 5585
$ws4.Range("F20").Value = 150
$ws4.Range("F22").Value = 995
$ws4.Range("F24").Value = 322
